# Update the "取得日時" (acquisition timestamp) column on the "ランサーズ"
# sheet for the newly (re-)appended rows: 2025-10-21 18:34:12 -> 2025-10-22 01:21:01
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-10-22 01:21:01"

for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
